$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.3327402135231317
$summary.Range("C2").Value = 0.06516290726817042
$summary.Range("D2").Value = 0.9285714285714286
$summary.Range("E2").Value = 0.1217798594847775
$summary.Range("F2").Value = 0.2544031311154599
$summary.Range("G2").Value = 0.6151046405823476
$summary.Range("H2").Value = 0.8103263777421081
$summary.Range("I2").Value = 26
$summary.Range("J2").Value = 373
$summary.Range("K2").Value = 161
$summary.Range("L2").Value = 2

# --- Sheet 2: Classification Report ---
$report = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$report.Range("B2").Value = 0.9877300613496932
$report.Range("C2").Value = 0.301498127340824
$report.Range("D2").Value = 0.4619799139167862

# Row 3 ("1")
$report.Range("B3").Value = 0.06516290726817042
$report.Range("C3").Value = 0.9285714285714286
$report.Range("D3").Value = 0.1217798594847775

# Row 4 ("accuracy")
$report.Range("B4").Value = 0.3327402135231317
$report.Range("C4").Value = 0.3327402135231317
$report.Range("D4").Value = 0.3327402135231317
$report.Range("E4").Value = 0.3327402135231317

# Row 5 ("macro avg")
$report.Range("B5").Value = 0.5264464843089318
$report.Range("C5").Value = 0.6150347779561263
$report.Range("D5").Value = 0.2918798867007819

# Row 6 ("weighted avg")
$report.Range("B6").Value = 0.9417658615022153
$report.Range("C6").Value = 0.3327402135231317
$report.Range("D6").Value = 0.4450304450127003

# --- Sheet 3: Confusion Matrix ---
$confusion = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$confusion.Range("B2").Value = 161
$confusion.Range("C2").Value = 373

# Row 3 ("Actual 1")
$confusion.Range("B3").Value = 2
$confusion.Range("C3").Value = 26
